$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 6565.833333333333, 7112, 6140, 0.1584909041722616),
    @(1, 6427.966666666666, 6829, 5647, 0.1591569185256958),
    @(2, 6250.4,            6755, 5387, 0.1612710078557332),
    @(3, 6833.4,            7438, 6475, 0.1599444309870402),
    @(4, 6334.6,            6855, 5532, 0.1662849505742391),
    @(5, 6318.966666666666, 6790, 5859, 0.1606734275817871),
    @(6, 6678.166666666667, 7159, 5812, 0.1638215621312459),
    @(7, 6623.933333333333, 7157, 5969, 0.1605967044830322),
    @(8, 6418.9,            6988, 5754, 0.1605101188023885),
    @(9, 6653.5,            7311, 5966, 0.1600111166636149)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
}
